{"js": "// Rewrite the bulleted \"Impact\" achievements under the \"KEY ACHIEVEMENTS AND IMPACT\"\n// heading so they read as short, impact-focused accomplishment statements instead of\n// the longer job-duty style bullets, and trim from 6 bullets down to 4.\n//\n// Old bullets (in order), now replaced:\n//  1) \"Built real-time FEC analysis systems ...\"              -> new text\n//  2) \"Built cloud-based data warehouse solutions on AWS ...\" -> new text\n//  3) \"Designed ETL pipelines using PySpark, dbt, ...\"        -> new text\n//  4) \"Trigonometric algorithm for boundary estimation ...\"   -> REMOVED\n//  5) \"Built redistricting platform used by thousands ...\"    -> REMOVED\n//  6) \"Discovered systematic race coding errors ...\"          -> new text\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Find the \"Impact\" sub-heading (Heading 3) that sits right under the\n// \"KEY ACHIEVEMENTS AND IMPACT\" section heading. Matching on style + exact\n// text keeps this from colliding with any other paragraph in the document.\nlet impactHeadingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Impact\" && p.style === \"Heading 3\") {\n    impactHeadingIndex = i;\n    break;\n  }\n}\n\nif (impactHeadingIndex === -1) {\n  throw new Error('Could not locate the \"Impact\" heading under KEY ACHIEVEMENTS AND IMPACT.');\n}\n\n// The six bullet paragraphs immediately follow that heading.\nconst bulletStart = impactHeadingIndex + 1;\nconst oldBullets = [\n  \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n];\n\n// Sanity-check that the paragraphs right after the heading match what we expect\n// before mutating anything.\nfor (let i = 0; i < oldBullets.length; i++) {\n  const actual = paragraphs.items[bulletStart + i].text;\n  if (actual !== oldBullets[i]) {\n    throw new Error(\n      \"Unexpected bullet text at position \" + i + \": \" + JSON.stringify(actual)\n    );\n  }\n}\n\n// New, shortened set of accomplishment-style bullets.\nconst newBullets = [\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\u2022 $4.7M savings enabled nonprofit access\",\n  \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"\u2022 178% accuracy improvement in racial classification algorithms\",\n];\n\n// Bullets 1-3 map 1:1 onto new bullets 1-3.\nparagraphs.items[bulletStart + 0].insertText(newBullets[0], \"Replace\");\nparagraphs.items[bulletStart + 1].insertText(newBullets[1], \"Replace\");\nparagraphs.items[bulletStart + 2].insertText(newBullets[2], \"Replace\");\n\n// Old bullet 6 becomes new bullet 4.\nparagraphs.items[bulletStart + 5].insertText(newBullets[3], \"Replace\");\n\n// Old bullets 4 and 5 are dropped entirely.\nparagraphs.items[bulletStart + 4].delete();\nparagraphs.items[bulletStart + 3].delete();\n\nawait context.sync();\n", "ps1": "# Rewrite the bulleted \"Impact\" achievements under the \"KEY ACHIEVEMENTS AND IMPACT\"\n# heading so they read as short, impact-focused accomplishment statements instead of\n# the longer job-duty style bullets, and trim from 6 bullets down to 4.\n#\n# Old bullets (in order), now replaced:\n#  1) \"Built real-time FEC analysis systems ...\"              -> new text\n#  2) \"Built cloud-based data warehouse solutions on AWS ...\" -> new text\n#  3) \"Designed ETL pipelines using PySpark, dbt, ...\"        -> new text\n#  4) \"Trigonometric algorithm for boundary estimation ...\"   -> REMOVED\n#  5) \"Built redistricting platform used by thousands ...\"    -> REMOVED\n#  6) \"Discovered systematic race coding errors ...\"          -> new text\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n\n# Locate the \"Impact\" sub-heading (Heading 3) that sits right under the\n# \"KEY ACHIEVEMENTS AND IMPACT\" section heading. Matching on style + exact\n# text keeps this from colliding with any other paragraph in the document\n# (e.g. the word \"impact\" also appears inside unrelated body text).\n$count = $d.Paragraphs.Count\n$impactHeadingIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"Impact\" -and $p.Style.NameLocal -eq \"Heading 3\") {\n        $impactHeadingIndex = $i\n        break\n    }\n}\n\nif ($impactHeadingIndex -eq -1) {\n    throw \"Could not locate the 'Impact' heading under KEY ACHIEVEMENTS AND IMPACT.\"\n}\n\n# The six bullet paragraphs immediately follow that heading (1-based indices).\n$bulletStart = $impactHeadingIndex + 1\n\n$oldBullets = @(\n    \"$bullet Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\",\n    \"$bullet Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n    \"$bullet Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n    \"$bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n    \"$bullet Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"\n)\n\n# Sanity-check that the paragraphs right after the heading match what we expect\n# before mutating anything.\nfor ($k = 0; $k -lt $oldBullets.Count; $k++) {\n    $actual = $d.Paragraphs.Item($bulletStart + $k).Range.Text.TrimEnd([char]13, [char]7)\n    if ($actual -ne $oldBullets[$k]) {\n        throw \"Unexpected bullet text at position $k`: $actual\"\n    }\n}\n\n# New, shortened set of accomplishment-style bullets.\n$newBullets = @(\n    \"$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    \"$bullet `$4.7M savings enabled nonprofit access\",\n    \"$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n    \"$bullet 178% accuracy improvement in racial classification algorithms\"\n)\n\n# Bullets 1-3 map 1:1 onto new bullets 1-3.\n$d.Paragraphs.Item($bulletStart + 0).Range.Text = $newBullets[0]\n$d.Paragraphs.Item($bulletStart + 1).Range.Text = $newBullets[1]\n$d.Paragraphs.Item($bulletStart + 2).Range.Text = $newBullets[2]\n\n# Old bullet 6 becomes new bullet 4.\n$d.Paragraphs.Item($bulletStart + 5).Range.Text = $newBullets[3]\n\n# Old bullets 4 and 5 are dropped entirely (delete higher index first so the\n# lower index stays valid).\n$d.Paragraphs.Item($bulletStart + 4).Range.Delete()\n$d.Paragraphs.Item($bulletStart + 3).Range.Delete()\n"}
